# Auto-generated edit script
# Applies the cell-value changes from the commit diff to the
# Garuda_Profits workbook's per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1703.6316
$ws.Range("I40").Value = 1545
$ws.Range("J40").Value = 1745.9333
$ws.Range("K40").Value = 1545
$ws.Range("L40").Value = 1745.9333
$ws.Range("M40").Value = -1370
$ws.Range("N40").Value = -2095.9333
# Row 64
$ws.Range("H64").Value = 3413.8462
$ws.Range("I64").Value = 2780
$ws.Range("J64").Value = 3466.6667
$ws.Range("K64").Value = 2780
$ws.Range("L64").Value = 3466.6667
$ws.Range("M64").Value = -2532
$ws.Range("N64").Value = -3962.6667
# Row 67
$ws.Range("H67").Value = 3413.8462
$ws.Range("I67").Value = 2780
$ws.Range("J67").Value = 3466.6667
$ws.Range("K67").Value = 2780
$ws.Range("L67").Value = 3466.6667
$ws.Range("M67").Value = -1922
$ws.Range("N67").Value = -5182.6667
# Row 76
$ws.Range("H76").Value = 73957.29
$ws.Range("I76").Value = 85808.586
$ws.Range("J76").Value = 2849.5
$ws.Range("K76").Value = 85808.586
$ws.Range("L76").Value = 2849.5
$ws.Range("M76").Value = -85493.586
$ws.Range("N76").Value = -3479.5
# Row 79
$ws.Range("H79").Value = 73957.29
$ws.Range("I79").Value = 85808.586
$ws.Range("J79").Value = 2849.5
$ws.Range("K79").Value = 85808.586
$ws.Range("L79").Value = 2849.5
$ws.Range("M79").Value = -84716.586
# Row 80
$ws.Range("H80").Value = 915.6
$ws.Range("I80").Value = 561.6
$ws.Range("J80").Value = 1269.6
$ws.Range("K80").Value = 1684.8
$ws.Range("L80").Value = 3808.8
$ws.Range("M80").Value = -686.8000000000002
$ws.Range("N80").Value = -5804.799999999999
# Row 83
$ws.Range("H83").Value = 915.6
$ws.Range("I83").Value = 561.6
$ws.Range("J83").Value = 1269.6
$ws.Range("K83").Value = 5054.400000000001
$ws.Range("L83").Value = 11426.4
$ws.Range("M83").Value = -62.40000000000055
$ws.Range("N83").Value = -21410.4
# Row 103
$ws.Range("H103").Value = 518.25
$ws.Range("I103").Value = 449.42856
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 1348.28568
$ws.Range("L103").Value = 3000
$ws.Range("M103").Value = -762.28568
$ws.Range("N103").Value = -4172

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 27027828
$ws.Range("I2").Value = 43478790
$ws.Range("J2").Value = 1249.1428
$ws.Range("K2").Value = 43478790
$ws.Range("L2").Value = 1249.1428
$ws.Range("M2").Value = -43478677
$ws.Range("N2").Value = -1475.1428
# Row 63
$ws.Range("H63").Value = 2000980
$ws.Range("I63").Value = 2500850
$ws.Range("K63").Value = 2500850
$ws.Range("M63").Value = -2500164
# Row 66
$ws.Range("H66").Value = 2000980
$ws.Range("I66").Value = 2500850
$ws.Range("K66").Value = 12504250
$ws.Range("M66").Value = -12500818
# Row 88
$ws.Range("H88").Value = 838641.5
$ws.Range("I88").Value = 1432612
$ws.Range("J88").Value = 7082.8
$ws.Range("K88").Value = 1432612
$ws.Range("L88").Value = 7082.8
$ws.Range("M88").Value = -1432206
$ws.Range("N88").Value = -7894.8
# Row 91
$ws.Range("H91").Value = 838641.5
$ws.Range("I91").Value = 1432612
$ws.Range("J91").Value = 7082.8
$ws.Range("K91").Value = 1432612
$ws.Range("L91").Value = 7082.8
$ws.Range("M91").Value = -1431208
$ws.Range("N91").Value = -9890.8
# Row 116
$ws.Range("H116").Value = 27027828
$ws.Range("I116").Value = 43478790
$ws.Range("J116").Value = 1249.1428
$ws.Range("K116").Value = 43478790
$ws.Range("L116").Value = 1249.1428
$ws.Range("M116").Value = -43476496
$ws.Range("N116").Value = -5837.1428
# Row 122
$ws.Range("H122").Value = 1915.3
$ws.Range("I122").Value = 1880.75
$ws.Range("J122").Value = 2053.5
$ws.Range("K122").Value = 5642.25
$ws.Range("L122").Value = 6160.5
$ws.Range("M122").Value = -3192.25
$ws.Range("N122").Value = -11060.5
# Row 132
$ws.Range("H132").Value = 5449.5674
$ws.Range("I132").Value = 5774.839
$ws.Range("K132").Value = 17324.517
$ws.Range("M132").Value = -14794.517

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 27027828
$ws.Range("I3").Value = 43478790
$ws.Range("J3").Value = 1249.1428
$ws.Range("K3").Value = 43478790
$ws.Range("L3").Value = 1249.1428
$ws.Range("M3").Value = -43478676
$ws.Range("N3").Value = -1477.1428

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 102
$ws.Range("H102").Value = 36455.555
$ws.Range("J102").Value = 36455.555
$ws.Range("L102").Value = 36455.555
$ws.Range("N102").Value = -41323.555

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 44415628
$ws.Range("I70").Value = 60403664
$ws.Range("J70").Value = 4412.1113
$ws.Range("K70").Value = 60403664
$ws.Range("L70").Value = 4412.1113
$ws.Range("M70").Value = -60403394
$ws.Range("N70").Value = -4952.1113
# Row 73
$ws.Range("H73").Value = 44415628
$ws.Range("I73").Value = 60403664
$ws.Range("J73").Value = 4412.1113
$ws.Range("K73").Value = 60403664
$ws.Range("L73").Value = 4412.1113
$ws.Range("M73").Value = -60402728
$ws.Range("N73").Value = -6284.1113
# Row 80
$ws.Range("H80").Value = 3711.6667
$ws.Range("I80").Value = 3711.6667
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3711.6667
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2713.6667
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 3711.6667
$ws.Range("I83").Value = 3711.6667
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 18558.3335
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -13566.3335
$ws.Range("N83").ClearContents()
# Row 122
$ws.Range("H122").Value = 50003604
$ws.Range("I122").Value = 111115510
$ws.Range("J122").Value = 2954.5454
$ws.Range("K122").Value = 333346530
$ws.Range("L122").Value = 8863.6362
$ws.Range("M122").Value = -333344080
$ws.Range("N122").Value = -13763.6362
# Row 123
$ws.Range("H123").Value = 19992.363
$ws.Range("J123").Value = 19992.363
$ws.Range("L123").Value = 19992.363
$ws.Range("N123").Value = -24892.363
# Row 131
$ws.Range("H131").Value = 21000
$ws.Range("J131").Value = 21000
$ws.Range("L131").Value = 21000
$ws.Range("N131").Value = -31080
# Row 132
$ws.Range("H132").Value = 254003
$ws.Range("I132").Value = 335337.34
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 1006012.02
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -1003482.02
$ws.Range("N132").Value = -35060

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 6133.0244
$ws.Range("I132").Value = 8563.885
$ws.Range("J132").Value = 1919.5333
$ws.Range("K132").Value = 25691.655
$ws.Range("L132").Value = 5758.5999
$ws.Range("M132").Value = -23161.655
$ws.Range("N132").Value = -10818.5999
# Row 138
$ws.Range("H138").Value = 54033.332
$ws.Range("J138").Value = 54033.332
$ws.Range("L138").Value = 54033.332
$ws.Range("N138").Value = -64313.332

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 106
$ws.Range("H106").Value = 26450
$ws.Range("J106").Value = 26450
$ws.Range("L106").Value = 26450
$ws.Range("N106").Value = -28974

